# Add "Planilha2" worksheet with test data mass for the search (pesquisa) tests
# (massa de dados for the magnifying-glass / home-page search test cases).
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Planilha2"

# Category column (A4:A7) first
$ws2.Cells.Item(4, 1).Value = "mice"
$ws2.Cells.Item(5, 1).Value = "headphones"
$ws2.Cells.Item(6, 1).Value = "tablets"
$ws2.Cells.Item(7, 1).Value = "speakers"

# Product-name column (B2, B4:B7)
$ws2.Cells.Item(2, 2).Value = "HP ENVY - 17T TOUCH LAPTOP"
$ws2.Cells.Item(4, 2).Value = "Logitech USB Headset H390"
$ws2.Cells.Item(5, 2).Value = "HP Elite x2 1011 G1 Tablet"
$ws2.Cells.Item(6, 2).Value = "HP Elite x2 1011 G1 Tablet"
$ws2.Cells.Item(7, 2).Value = "HP Roar Plus Wireless Speaker"

# Row 2 - laptops example (category + id)
$ws2.Cells.Item(2, 1).Value = "laptops"
$ws2.Cells.Item(2, 3).Value = 7

# Header row
$ws2.Cells.Item(1, 1).Value = "pesquisa"
$ws2.Cells.Item(1, 2).Value = "nome produto"
$ws2.Cells.Item(1, 3).Value = "id"

# Leftover inline formatting from the pasted web data (dark-gray text,
# some cells additionally underlined) -- matches what a browser->Excel
# paste of the product table leaves behind.
$ws2.Cells.Item(1, 3).Font.Color = 2236962
$ws2.Cells.Item(2, 4).Font.Color = 2236962
$ws2.Cells.Item(5, 1).Font.Color = 2236962
$ws2.Cells.Item(6, 2).Font.Color = 2236962

$ws2.Cells.Item(2, 5).Font.Color = 2236962
$ws2.Cells.Item(2, 5).Font.Underline = 2
$ws2.Cells.Item(7, 2).Font.Color = 2236962
$ws2.Cells.Item(7, 2).Font.Underline = 2

# Column widths (bestFit) to match target layout
$ws2.Columns.Item(1).ColumnWidth = 11.453125
$ws2.Columns.Item(2).ColumnWidth = 27
$ws2.Columns.Item(3).ColumnWidth = 22.54296875
$ws2.Columns.Item(4).ColumnWidth = 22.54296875
$ws2.Columns.Item(5).ColumnWidth = 26.26953125

# Move the new sheet so it sits after Planilha1 (do this AFTER writing all
# data -- the COM object handle otherwise can resolve by stale index).
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2.Move($null, $ws1)

# Planilha1 keeps its own pre-existing selection. Apply this BEFORE activating
# Planilha2 below -- selecting a range on a sheet implicitly activates that
# sheet/window, so doing it last would steal the active tab back.
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws1.Range("E5").Select()

# Select + activate Planilha2 last so it ends up the active/displayed tab.
$ws2 = $wb.Worksheets.Item("Planilha2")
$ws2.Range("C2").Select()
$ws2.Activate()
